$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 17.333334
$ws.Range("H42").Value = 518.9286
$ws.Range("I42").Value = 122.42857
$ws.Range("K42").Value = 367.28571
$ws.Range("M42").Value = -137.28571
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -40924
$ws.Range("H98").Value = 1898.4546
$ws.Range("I98").Value = 1930.8235
$ws.Range("J98").Value = 1788.4
$ws.Range("K98").Value = 1930.8235
$ws.Range("L98").Value = 1788.4
$ws.Range("M98").Value = -432.8235
$ws.Range("N98").Value = -4784.4
$ws.Range("H112").Value = 2952.6538
$ws.Range("I112").Value = 1624.5
$ws.Range("J112").Value = 3005.78
$ws.Range("K112").Value = 4873.5
$ws.Range("L112").Value = 9017.34
$ws.Range("M112").Value = -3765.5
$ws.Range("N112").Value = -11233.34
$ws.Range("H122").Value = 1898.4546
$ws.Range("I122").Value = 1930.8235
$ws.Range("J122").Value = 1788.4
$ws.Range("K122").Value = 5792.470499999999
$ws.Range("L122").Value = 5365.200000000001
$ws.Range("M122").Value = -3342.470499999999
$ws.Range("N122").Value = -10265.2
$ws.Range("H138").Value = 6532.915
$ws.Range("I138").Value = 2879.375
$ws.Range("K138").Value = 8638.125
$ws.Range("M138").Value = -3498.125

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 34022.4
$ws.Range("I88").Value = 698.4286
$ws.Range("J88").Value = 63180.875
$ws.Range("K88").Value = 698.4286
$ws.Range("L88").Value = 63180.875
$ws.Range("M88").Value = -292.4286
$ws.Range("N88").Value = -63992.875
$ws.Range("H91").Value = 34022.4
$ws.Range("I91").Value = 698.4286
$ws.Range("J91").Value = 63180.875
$ws.Range("K91").Value = 698.4286
$ws.Range("L91").Value = 63180.875
$ws.Range("M91").Value = 705.5714
$ws.Range("N91").Value = -65988.875
$ws.Range("H96").Value = 60000
$ws.Range("J96").Value = 60000
$ws.Range("L96").Value = 60000
$ws.Range("N96").Value = -65492
$ws.Range("H132").Value = 5426.5186
$ws.Range("I132").Value = 2261.2144
$ws.Range("K132").Value = 6783.6432
$ws.Range("M132").Value = -4253.6432
$ws.Range("H140").Value = 110165.5
$ws.Range("J140").Value = 110165.5
$ws.Range("L140").Value = 110165.5
$ws.Range("N140").Value = -120525.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 104988.2
$ws.Range("J132").Value = 104988.2
$ws.Range("L132").Value = 104988.2
$ws.Range("N132").Value = -115108.2
$ws.Range("H134").Value = 2925.4
$ws.Range("I134").Value = 2777.2144
$ws.Range("K134").Value = 8331.643199999999
$ws.Range("M134").Value = -5796.643199999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3617.5898
$ws.Range("I31").Value = 3163.0984
$ws.Range("J31").Value = 5248.4116
$ws.Range("K31").Value = 3163.0984
$ws.Range("L31").Value = 5248.4116
$ws.Range("M31").Value = -2868.0984
$ws.Range("N31").Value = -5838.4116
$ws.Range("H34").Value = 3617.5898
$ws.Range("I34").Value = 3163.0984
$ws.Range("J34").Value = 5248.4116
$ws.Range("K34").Value = 3163.0984
$ws.Range("L34").Value = 5248.4116
$ws.Range("M34").Value = -2961.0984
$ws.Range("N34").Value = -5652.4116
$ws.Range("H93").Value = 65316.668
$ws.Range("I93").Value = 6000
$ws.Range("J93").Value = 77180
$ws.Range("K93").Value = 6000
$ws.Range("L93").Value = 77180
$ws.Range("M93").Value = -4128
$ws.Range("N93").Value = -80924

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2149.8333
$ws.Range("I25").Value = 1750
$ws.Range("J25").Value = 2349.75
$ws.Range("K25").Value = 5250
$ws.Range("L25").Value = 7049.25
$ws.Range("M25").Value = -5081
$ws.Range("N25").Value = -7387.25
$ws.Range("H26").Value = 850.4
$ws.Range("I26").Value = 313
$ws.Range("K26").Value = 939
$ws.Range("M26").Value = -651
$ws.Range("H30").Value = 2149.8333
$ws.Range("I30").Value = 1750
$ws.Range("J30").Value = 2349.75
$ws.Range("K30").Value = 5250
$ws.Range("L30").Value = 7049.25
$ws.Range("M30").Value = -5148
$ws.Range("N30").Value = -7253.25
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("H57").Value = 14794.5
$ws.Range("I57").Value = 4592
$ws.Range("J57").Value = 24997
$ws.Range("K57").Value = 13776
$ws.Range("L57").Value = 74991
$ws.Range("M57").Value = -13217
$ws.Range("N57").Value = -76109
$ws.Range("H62").Value = 19991
$ws.Range("J62").Value = 19991
$ws.Range("L62").Value = 59973
$ws.Range("N62").Value = -61345
$ws.Range("H65").Value = 19991
$ws.Range("J65").Value = 19991
$ws.Range("L65").Value = 179919
$ws.Range("N65").Value = -186783
$ws.Range("H96").Value = 12998.833
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 12998.833
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 38996.499
$ws.Range("N96").Value = -43114.499
$ws.Range("N32").ClearContents()
$ws.Range("M96").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15050000
$ws.Range("I11").Value = 13400000
$ws.Range("K11").Value = 13400000
$ws.Range("M11").Value = -13399861
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("H113").Value = 13023.5
$ws.Range("I113").Value = 11663
$ws.Range("K113").Value = 11663
$ws.Range("M113").Value = -9493
$ws.Range("H122").Value = 919484.5
$ws.Range("I122").Value = 1002801.25
$ws.Range("K122").Value = 3008403.75
$ws.Range("M122").Value = -3005953.75
$ws.Range("H126").Value = 3549.7742
$ws.Range("J126").Value = 6076.636
$ws.Range("L126").Value = 18229.908
$ws.Range("N126").Value = -23169.908
$ws.Range("H132").Value = 7542.7393
$ws.Range("I132").Value = 6912.2
$ws.Range("J132").Value = 8725
$ws.Range("K132").Value = 20736.6
$ws.Range("L132").Value = 26175
$ws.Range("M132").Value = -18206.6
$ws.Range("N52").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1191.9286
$ws.Range("I22").Value = 817.8
$ws.Range("K22").Value = 817.8
$ws.Range("M22").Value = -522.8
$ws.Range("H27").Value = 1191.9286
$ws.Range("I27").Value = 817.8
$ws.Range("K27").Value = 817.8
$ws.Range("M27").Value = -710.8
$ws.Range("H132").Value = 5609.485
$ws.Range("I132").Value = 3799.125
$ws.Range("J132").Value = 7313.353
$ws.Range("K132").Value = 11397.375
$ws.Range("L132").Value = 21940.059
$ws.Range("M132").Value = -8867.375
$ws.Range("N132").Value = -27000.059

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28593
$ws.Range("J41").Value = 28593
$ws.Range("L41").Value = 28593
$ws.Range("N41").Value = -29373
$ws.Range("H81").Value = 3476291
$ws.Range("I81").Value = 2609228.5
$ws.Range("J81").Value = 5210416
$ws.Range("K81").Value = 5218457
$ws.Range("L81").Value = 10420832
$ws.Range("M81").Value = -5217396
$ws.Range("N81").Value = -10422954
$ws.Range("H84").Value = 3476291
$ws.Range("I84").Value = 2609228.5
$ws.Range("J84").Value = 5210416
$ws.Range("K84").Value = 26092285
$ws.Range("L84").Value = 52104160
$ws.Range("M84").Value = -26086981
$ws.Range("N84").Value = -52114768
$ws.Range("H132").Value = 12505207
$ws.Range("I132").Value = 1658.2812
$ws.Range("J132").Value = 62519404
$ws.Range("K132").Value = 4974.8436
$ws.Range("L132").Value = 187558212
$ws.Range("M132").Value = -2444.8436
$ws.Range("N132").Value = -187563272
